$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add the timesheet entries for 31 January 2020.
# Layout mirrors the existing day blocks: a blank separator row (formatted
# like the header-ish separator rows, e.g. row 20), followed by one row per
# hour worked with Timestamp / Task / Location columns.
# ---------------------------------------------------------------------------

# Row 29: blank separator row (same look as row 20 -> styles 5/6/5)
$ws.Range("A20:C20").Copy()
$ws.Range("A29").PasteSpecial(-4122)

# Row 30: 10:00 - 11:00 (wrapped task text like row 19 -> styles 1/2/1)
$ws.Range("A19:C19").Copy()
$ws.Range("A30").PasteSpecial(-4122)
$ws.Range("A30").Value = "Jan 31 10:00 to 11:00"
$ws.Range("B30").Value = "Trying to convert buckets to dictionaries but the format of buckets
sample is conflicted while writing to dict. No json or ast methods are
applicable to it. Working on another solution to re-create it."
$ws.Range("C30").Value = "Infimetrics"
$ws.Rows.Item(30).RowHeight = 60

# Row 31: 11:00 - 12:00 (wrapped task text -> styles 1/2/1)
$ws.Range("A19:C19").Copy()
$ws.Range("A31").PasteSpecial(-4122)
$ws.Range("A31").Value = "Jan 31 11:00 to 12:00"
$ws.Range("B31").Value = "Figured out solution of loading buckets as dictionaries by transforming
buckets to list from data exploration of sysytem health file. Using json
loads to make bucket into dictionary. For this modified code creation 
of bucket data by writing starting and ending time as string."
$ws.Range("C31").Value = "Infimetrics"
$ws.Rows.Item(31).RowHeight = 105

# Row 32: 12:00 - 13:00 (single line task text like row 26 -> styles 1/3/1)
$ws.Range("A26:C26").Copy()
$ws.Range("A32").PasteSpecial(-4122)
$ws.Range("A32").Value = "Jan 31 12:00 to 13:00"
$ws.Range("B32").Value = "Creating percentile buckets"
$ws.Range("C32").Value = "Infimetrics"

# Row 33: 13:00 - 14:00 (wrapped task text -> styles 1/2/1)
$ws.Range("A19:C19").Copy()
$ws.Range("A33").PasteSpecial(-4122)
$ws.Range("A33").Value = "Jan 31 13:00 to 14:00"
$ws.Range("B33").Value = "Created percentile buckets, choosed percentile as 80, facing some 
issues regarding of nan values. "
$ws.Range("C33").Value = "Infimetrics"
$ws.Rows.Item(33).RowHeight = 30

# Row 34: 14:00 - 15:00 (single line "Lunch" task like row 3 -> styles 1/3/1)
$ws.Range("A26:C26").Copy()
$ws.Range("A34").PasteSpecial(-4122)
$ws.Range("A34").Value = "Jan 31 14:00 to 15:00"
$ws.Range("B34").Value = "Lunch"
$ws.Range("C34").Value = "Infimetrics"

$excel.CutCopyMode = $false

# Update the view: select C35 so the saved selection state recorded in the
# workbook matches (mirrors where the author's cursor ended up).
$null = $ws.Range("C35").Select()
